$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The ENTSO-E consumption extract has rolled forward by 4 days:
#   12.06.2025 / 13.06.2025  ->  16.06.2025 / 17.06.2025
# This updates:
#   - Column D ("Lookup") shared-string labels (date + quarter-index)
#   - Column A ("Timestamp") Excel serial date/time values (+4 days)
#   - Column B ("Actual Consumption (MW)") with the newly fetched figures
# ---------------------------------------------------------------------------

# --- Column D: Lookup labels, day 1 (rows 2-97) and day 2 (rows 98-193) ---
For ($i = 2; $i -le 97; $i++) {
    $q = $i - 1
    $ws.Cells.Item($i, 4).Value = "16.06.2025" + $q
}
For ($i = 98; $i -le 193; $i++) {
    $q = $i - 97
    $ws.Cells.Item($i, 4).Value = "17.06.2025" + $q
}

# --- Columns A (Timestamp) and B (Actual Consumption (MW)) ---
$aValues = @(45824,45824.01041666666,45824.02083333334,45824.03125,45824.04166666666,45824.05208333334,45824.0625,45824.07291666666,45824.08333333334,45824.09375,45824.10416666666,45824.11458333334,45824.125,45824.13541666666,45824.14583333334,45824.15625,45824.16666666666,45824.17708333334,45824.1875,45824.19791666666,45824.20833333334,45824.21875,45824.22916666666,45824.23958333334,45824.25,45824.26041666666,45824.27083333334,45824.28125,45824.29166666666,45824.30208333334,45824.3125,45824.32291666666,45824.33333333334,45824.34375,45824.35416666666,45824.36458333334,45824.375,45824.38541666666,45824.39583333334,45824.40625,45824.41666666666,45824.42708333334,45824.4375,45824.44791666666,45824.45833333334,45824.46875,45824.47916666666,45824.48958333334,45824.5,45824.51041666666,45824.52083333334,45824.53125,45824.54166666666,45824.55208333334,45824.5625,45824.57291666666,45824.58333333334,45824.59375,45824.60416666666,45824.61458333334,45824.625,45824.63541666666,45824.64583333334,45824.65625,45824.66666666666,45824.67708333334,45824.6875,45824.69791666666,45824.70833333334,45824.71875,45824.72916666666,45824.73958333334,45824.75,45824.76041666666,45824.77083333334,45824.78125,45824.79166666666,45824.80208333334,45824.8125,45824.82291666666,45824.83333333334,45824.84375,45824.85416666666,45824.86458333334,45824.875,45824.88541666666,45824.89583333334,45824.90625,45824.91666666666,45824.92708333334,45824.9375,45824.94791666666,45824.95833333334,45824.96875,45824.97916666666,45824.98958333334,45825,45825.01041666666,45825.02083333334,45825.03125,45825.04166666666,45825.05208333334,45825.0625,45825.07291666666,45825.08333333334,45825.09375,45825.10416666666,45825.11458333334,45825.125,45825.13541666666,45825.14583333334,45825.15625,45825.16666666666,45825.17708333334,45825.1875,45825.19791666666,45825.20833333334,45825.21875,45825.22916666666,45825.23958333334,45825.25,45825.26041666666,45825.27083333334,45825.28125,45825.29166666666,45825.30208333334,45825.3125,45825.32291666666,45825.33333333334,45825.34375,45825.35416666666,45825.36458333334,45825.375,45825.38541666666,45825.39583333334,45825.40625,45825.41666666666,45825.42708333334,45825.4375,45825.44791666666,45825.45833333334,45825.46875,45825.47916666666,45825.48958333334,45825.5,45825.51041666666,45825.52083333334,45825.53125,45825.54166666666,45825.55208333334,45825.5625,45825.57291666666,45825.58333333334,45825.59375,45825.60416666666,45825.61458333334,45825.625,45825.63541666666,45825.64583333334,45825.65625,45825.66666666666,45825.67708333334,45825.6875,45825.69791666666,45825.70833333334,45825.71875,45825.72916666666,45825.73958333334,45825.75,45825.76041666666,45825.77083333334,45825.78125,45825.79166666666,45825.80208333334,45825.8125,45825.82291666666,45825.83333333334,45825.84375,45825.85416666666,45825.86458333334,45825.875,45825.88541666666,45825.89583333334,45825.90625,45825.91666666666,45825.92708333334,45825.9375,45825.94791666666,45825.95833333334,45825.96875,45825.97916666666,45825.98958333334)
$bValues = @(4813,4710,4669,4624,4598,4595,4569,4591,4565,4551,4551,4555,4532,4555,4602,4618,4637,4668,4669,4736,4895,5057,5181,5295,5532,5573,5626,5652,5752,5738,5751,5694,5490,5483,5471,5448,5335,5251,5227,5236,5040,5050,4999,5036,4905,4943,4874,4949,4945,5031,4982,4950,4990,5026,4983,5042,5134,5227,5260,5269,5364,5472,5538,5546,5686,5764,5827,5978,6080,6203,6355,6440,6562,6727,6860,6998,6983,7031,7024,7092,7145,7164,7198,7183,7096,6983,6727,6551,6355,6273,6132,5978,5810,5759,5656,5582,5509,5420,5377,5365,5329,5315,5301,5299,5237,5235,5218,5225,5309,5257,5343,5300,5262,5271,5258,5357,5451,5528,5604,5658,5857,5932,5960,5948,6026,5993,5982,5982,5928,5849,5862,5802,5688,5634,5597,5557,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

For ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
